$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy header style (bold, bordered, centered) from existing header cell H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats

# Data for columns I (I0) and J (IF), rows 2-42
$values = @(
    @(7, 8),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(4, 5),
    @(7, 7),
    @(8, 8),
    @(7, 8),
    @(8, 9),
    @(7, 7),
    @(6, 7),
    @(7, 8),
    @(8, 8),
    @(7, 7),
    @(4, 5),
    @(7, 7),
    @(7, 8),
    @(6, 7),
    @(6, 6),
    @(8, 8),
    @(6, 7),
    @(8, 8),
    @(6, 6),
    @(9, 10),
    @(8, 9),
    @(5, 6),
    @(9, 9),
    @(6, 6),
    @(9, 9),
    @(7, 8),
    @(7, 8),
    @(8, 8),
    @(8, 9),
    @(6, 7),
    @(6, 6),
    @(6, 7),
    @(2, 3),
    @(3, 3)
)

$row = 2
foreach ($pair in $values) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
